$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# Row 17 (ALC)
$ws_ALC.Range("H17").Value = 1121.8718
$ws_ALC.Range("J17").Value = 860.2632
$ws_ALC.Range("L17").Value = 2580.7896
$ws_ALC.Range("N17").Value = -2916.7896

# Row 121 (ALC)
$ws_ALC.Range("H121").Value = 999
$ws_ALC.Range("J121").Value = 999
$ws_ALC.Range("L121").Value = 2997
$ws_ALC.Range("N121").Value = -6491

# Row 138 (ALC)
$ws_ALC.Range("H138").Value = 2432.4468
$ws_ALC.Range("I138").Value = 2027.7
$ws_ALC.Range("J138").Value = 3146.7058
$ws_ALC.Range("K138").Value = 6083.1
$ws_ALC.Range("L138").Value = 9440.117400000001
$ws_ALC.Range("M138").Value = -943.1000000000004
$ws_ALC.Range("N138").Value = -19720.1174

$ws_ARM = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$ws_ARM.Range("H2").Value = 5814677
$ws_ARM.Range("J2").Value = 999
$ws_ARM.Range("L2").Value = 999
$ws_ARM.Range("N2").Value = -1225

# Row 4 (ARM)
$ws_ARM.Range("H4").Value = 27.5
$ws_ARM.Range("I4").Value = 27.5
$ws_ARM.Range("K4").Value = 27.5
$ws_ARM.Range("M4").Value = 88.5

# Row 23 (ARM)
$ws_ARM.Range("H23").Value = 26950.834
$ws_ARM.Range("I23").Value = 45003
$ws_ARM.Range("J23").Value = 17924.75
$ws_ARM.Range("K23").Value = 45003
$ws_ARM.Range("L23").Value = 17924.75
$ws_ARM.Range("M23").Value = -44744
$ws_ARM.Range("N23").Value = -18442.75

# Row 30 (ARM)
$ws_ARM.Range("M30").Value = -650
$ws_ARM.Range("H30").Value = 800
$ws_ARM.Range("I30").Value = 800
$ws_ARM.Range("K30").Value = 800

# Row 37 (ARM)
$ws_ARM.Range("H37").Value = 14032.571
$ws_ARM.Range("J37").Value = 14032.571
$ws_ARM.Range("L37").Value = 14032.571
$ws_ARM.Range("N37").Value = -14578.571

# Row 74 (ARM)
$ws_ARM.Range("H74").Value = 940.87805
$ws_ARM.Range("I74").Value = 768.75757
$ws_ARM.Range("J74").Value = 1650.875
$ws_ARM.Range("K74").Value = 768.75757
$ws_ARM.Range("L74").Value = 1650.875
$ws_ARM.Range("M74").Value = 105.24243
$ws_ARM.Range("N74").Value = -3398.875

# Row 77 (ARM)
$ws_ARM.Range("H77").Value = 940.87805
$ws_ARM.Range("I77").Value = 768.75757
$ws_ARM.Range("J77").Value = 1650.875
$ws_ARM.Range("K77").Value = 3843.78785
$ws_ARM.Range("L77").Value = 8254.375
$ws_ARM.Range("M77").Value = 524.2121500000003
$ws_ARM.Range("N77").Value = -16990.375

# Row 110 (ARM)
$ws_ARM.Range("H110").Value = 1930.5294
$ws_ARM.Range("I110").Value = 1593
$ws_ARM.Range("J110").Value = 2412.7144
$ws_ARM.Range("K110").Value = 1593
$ws_ARM.Range("L110").Value = 2412.7144
$ws_ARM.Range("M110").Value = 452
$ws_ARM.Range("N110").Value = -6502.7144

# Row 116 (ARM)
$ws_ARM.Range("H116").Value = 5814677
$ws_ARM.Range("J116").Value = 999
$ws_ARM.Range("L116").Value = 999
$ws_ARM.Range("N116").Value = -5587

# Row 132 (ARM)
$ws_ARM.Range("H132").Value = 1445.475
$ws_ARM.Range("I132").Value = 973.05
$ws_ARM.Range("J132").Value = 1917.9
$ws_ARM.Range("K132").Value = 2919.15
$ws_ARM.Range("L132").Value = 5753.700000000001
$ws_ARM.Range("M132").Value = -389.1499999999996
$ws_ARM.Range("N132").Value = -10813.7

$ws_BSM = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$ws_BSM.Range("H3").Value = 5814677
$ws_BSM.Range("J3").Value = 999
$ws_BSM.Range("L3").Value = 999
$ws_BSM.Range("N3").Value = -1227

# Row 64 (BSM)
$ws_BSM.Range("H64").Value = 712.2857
$ws_BSM.Range("I64").Value = 771.5
$ws_BSM.Range("J64").Value = 633.3333
$ws_BSM.Range("K64").Value = 771.5
$ws_BSM.Range("L64").Value = 633.3333
$ws_BSM.Range("M64").Value = -546.5
$ws_BSM.Range("N64").Value = -1083.3333

# Row 67 (BSM)
$ws_BSM.Range("H67").Value = 712.2857
$ws_BSM.Range("I67").Value = 771.5
$ws_BSM.Range("J67").Value = 633.3333
$ws_BSM.Range("K67").Value = 771.5
$ws_BSM.Range("L67").Value = 633.3333
$ws_BSM.Range("M67").Value = 8.5
$ws_BSM.Range("N67").Value = -2193.3333

# Row 99 (BSM)
$ws_BSM.Range("H99").Value = 1429.5714
$ws_BSM.Range("I99").Value = 951.75
$ws_BSM.Range("J99").Value = 2066.6667
$ws_BSM.Range("K99").Value = 951.75
$ws_BSM.Range("L99").Value = 2066.6667
$ws_BSM.Range("M99").Value = 546.25
$ws_BSM.Range("N99").Value = -5062.6667

$ws_CRP = $wb.Worksheets.Item("CRP")
# Row 8 (CRP)
$ws_CRP.Range("H8").Value = 3873.3333
$ws_CRP.Range("J8").Value = 3873.3333
$ws_CRP.Range("L8").Value = 3873.3333
$ws_CRP.Range("N8").Value = -4153.3333

# Row 11 (CRP)
$ws_CRP.Range("M11").ClearContents()
$ws_CRP.Range("H11").Value = 0
$ws_CRP.Range("I11").Value = 0
$ws_CRP.Range("K11").Value = 0

# Row 51 (CRP)
$ws_CRP.Range("H51").Value = 33000
$ws_CRP.Range("J51").Value = 33000
$ws_CRP.Range("L51").Value = 33000
$ws_CRP.Range("N51").Value = -34472

# Row 60 (CRP)
$ws_CRP.Range("H60").Value = 11333.444
$ws_CRP.Range("J60").Value = 11333.444
$ws_CRP.Range("L60").Value = 11333.444
$ws_CRP.Range("N60").Value = -12355.444

# Row 61 (CRP)
$ws_CRP.Range("H61").Value = 33000
$ws_CRP.Range("J61").Value = 33000
$ws_CRP.Range("L61").Value = 33000
$ws_CRP.Range("N61").Value = -33696

# Row 92 (CRP)
$ws_CRP.Range("H92").Value = 30495
$ws_CRP.Range("J92").Value = 30495
$ws_CRP.Range("L92").Value = 30495
$ws_CRP.Range("N92").Value = -35487

$ws_CUL = $wb.Worksheets.Item("CUL")
# Row 5 (CUL)
$ws_CUL.Range("H5").Value = 798.44446
$ws_CUL.Range("I5").Value = 713.25
$ws_CUL.Range("J5").Value = 866.6
$ws_CUL.Range("K5").Value = 2139.75
$ws_CUL.Range("L5").Value = 2599.8
$ws_CUL.Range("M5").Value = -2027.75
$ws_CUL.Range("N5").Value = -2823.8

# Row 131 (CUL)
$ws_CUL.Range("H131").Value = 9982.855
$ws_CUL.Range("I131").Value = 660.7143
$ws_CUL.Range("J131").Value = 10928.58
$ws_CUL.Range("K131").Value = 1982.1429
$ws_CUL.Range("L131").Value = 32785.74
$ws_CUL.Range("M131").Value = 3057.8571
$ws_CUL.Range("N131").Value = -42865.74

# Row 135 (CUL)
$ws_CUL.Range("H135").Value = 798.44446
$ws_CUL.Range("I135").Value = 713.25
$ws_CUL.Range("J135").Value = 866.6
$ws_CUL.Range("K135").Value = 6419.25
$ws_CUL.Range("L135").Value = 7799.400000000001
$ws_CUL.Range("M135").Value = -3884.25
$ws_CUL.Range("N135").Value = -12869.4

$ws_GSM = $wb.Worksheets.Item("GSM")
# Row 132 (GSM)
$ws_GSM.Range("H132").Value = 1243400.8
$ws_GSM.Range("I132").Value = 1674160.1
$ws_GSM.Range("J132").Value = 4967.5
$ws_GSM.Range("K132").Value = 5022480.300000001
$ws_GSM.Range("L132").Value = 14902.5
$ws_GSM.Range("M132").Value = -5019950.300000001
$ws_GSM.Range("N132").Value = -19962.5

$ws_LTW = $wb.Worksheets.Item("LTW")
# Row 2 (LTW)
$ws_LTW.Range("H2").Value = 243431.81
$ws_LTW.Range("J2").Value = 29625
$ws_LTW.Range("L2").Value = 29625
$ws_LTW.Range("N2").Value = -29849

# Row 22 (LTW)
$ws_LTW.Range("H22").Value = 2043.5385
$ws_LTW.Range("I22").Value = 5100
$ws_LTW.Range("J22").Value = 1487.8182
$ws_LTW.Range("K22").Value = 5100
$ws_LTW.Range("L22").Value = 1487.8182
$ws_LTW.Range("M22").Value = -4805
$ws_LTW.Range("N22").Value = -2077.8182

# Row 27 (LTW)
$ws_LTW.Range("H27").Value = 2043.5385
$ws_LTW.Range("I27").Value = 5100
$ws_LTW.Range("J27").Value = 1487.8182
$ws_LTW.Range("K27").Value = 5100
$ws_LTW.Range("L27").Value = 1487.8182
$ws_LTW.Range("M27").Value = -4993
$ws_LTW.Range("N27").Value = -1701.8182

$ws_WVR = $wb.Worksheets.Item("WVR")
# Row 82 (WVR)
$ws_WVR.Range("H82").Value = 70000
$ws_WVR.Range("J82").Value = 70000
$ws_WVR.Range("L82").Value = 70000
$ws_WVR.Range("N82").Value = -70766

# Row 85 (WVR)
$ws_WVR.Range("H85").Value = 70000
$ws_WVR.Range("J85").Value = 70000
$ws_WVR.Range("L85").Value = 70000
$ws_WVR.Range("N85").Value = -72652

# Row 113 (WVR)
$ws_WVR.Range("H113").Value = 1600
$ws_WVR.Range("I113").Value = 1200
$ws_WVR.Range("J113").Value = 2000
$ws_WVR.Range("K113").Value = 3600
$ws_WVR.Range("L113").Value = 6000
$ws_WVR.Range("M113").Value = -1430
$ws_WVR.Range("N113").Value = -10340

# Row 132 (WVR)
$ws_WVR.Range("H132").Value = 1790
$ws_WVR.Range("I132").Value = 1388.08
$ws_WVR.Range("J132").Value = 3799.6
$ws_WVR.Range("K132").Value = 4164.24
$ws_WVR.Range("L132").Value = 11398.8
$ws_WVR.Range("M132").Value = -1634.24
$ws_WVR.Range("N132").Value = -16458.8
